$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two Bortle-class descriptions that were reworded to mention
# "the Triangulum Galaxy, M33" explicitly.
$ws.Range("C2").Value = "The zodiacal light, gegenschein, and zodiacal band  are all visible — the zodiacal light to a striking degree, and the zodiacal band spanning the entire sky. Even with direct vision, the Triangulum Galaxy, M33, is is an obvious naked-eye object. The Scorpius and Sagittarius region of the Milky Way casts obvious diffuse shadows on the ground."
$ws.Range("C3").Value = "Airglow may be weakly apparent along the horizon. The Triangulum Galaxy, M33, is easily seen with direct vision. The summer Milky Way is highly structured to the unaided eye, and its brightest parts look like veined marble when viewed with ordinary binoculars."

# Move the active selection, matching the cursor position left behind by the author.
$ws.Range("G19").Select()
